$d = $word.ActiveDocument

$d.Content.Find.Execute("64×20=1280", $true, $false, $false, $false, $false, $true, 1, $false, "67×37=2479", 2) | Out-Null
$d.Content.Find.Execute("61×67=4087", $true, $false, $false, $false, $false, $true, 1, $false, "57×56=3192", 2) | Out-Null
$d.Content.Find.Execute("95×73=6935", $true, $false, $false, $false, $false, $true, 1, $false, "93×69=6417", 2) | Out-Null
$d.Content.Find.Execute("26×88=2288", $true, $false, $false, $false, $false, $true, 1, $false, "56×54=3024", 2) | Out-Null
$d.Content.Find.Execute("42×82=3444", $true, $false, $false, $false, $false, $true, 1, $false, "88×21=1848", 2) | Out-Null
$d.Content.Find.Execute("99×24=2376", $true, $false, $false, $false, $false, $true, 1, $false, "67×69=4623", 2) | Out-Null
$d.Content.Find.Execute("84×18=1512", $true, $false, $false, $false, $false, $true, 1, $false, "69×16=1104", 2) | Out-Null
$d.Content.Find.Execute("68×42=2856", $true, $false, $false, $false, $false, $true, 1, $false, "38×31=1178", 2) | Out-Null
$d.Content.Find.Execute("62×80=4960", $true, $false, $false, $false, $false, $true, 1, $false, "74×85=6290", 2) | Out-Null
$d.Content.Find.Execute("28×27=756", $true, $false, $false, $false, $false, $true, 1, $false, "87×32=2784", 2) | Out-Null
$d.Content.Find.Execute("35×52=1820", $true, $false, $false, $false, $false, $true, 1, $false, "42×67=2814", 2) | Out-Null
$d.Content.Find.Execute("66×51=3366", $true, $false, $false, $false, $false, $true, 1, $false, "88×52=4576", 2) | Out-Null
$d.Content.Find.Execute("76×50=3800", $true, $false, $false, $false, $false, $true, 1, $false, "87×14=1218", 2) | Out-Null
$d.Content.Find.Execute("45×65=2925", $true, $false, $false, $false, $false, $true, 1, $false, "48×94=4512", 2) | Out-Null
$d.Content.Find.Execute("60×18=1080", $true, $false, $false, $false, $false, $true, 1, $false, "99×16=1584", 2) | Out-Null
$d.Content.Find.Execute("92×50=4600", $true, $false, $false, $false, $false, $true, 1, $false, "26×71=1846", 2) | Out-Null
$d.Content.Find.Execute("31×34=1054", $true, $false, $false, $false, $false, $true, 1, $false, "52×94=4888", 2) | Out-Null
$d.Content.Find.Execute("70×69=4830", $true, $false, $false, $false, $false, $true, 1, $false, "19×18=342", 2) | Out-Null
$d.Content.Find.Execute("94×40=3760", $true, $false, $false, $false, $false, $true, 1, $false, "35×71=2485", 2) | Out-Null
$d.Content.Find.Execute("76×56=4256", $true, $false, $false, $false, $false, $true, 1, $false, "74×64=4736", 2) | Out-Null
$d.Content.Find.Execute("35×73=2555", $true, $false, $false, $false, $false, $true, 1, $false, "51×59=3009", 2) | Out-Null
$d.Content.Find.Execute("62×32=1984", $true, $false, $false, $false, $false, $true, 1, $false, "28×98=2744", 2) | Out-Null
$d.Content.Find.Execute("35×26=910", $true, $false, $false, $false, $false, $true, 1, $false, "60×95=5700", 2) | Out-Null
$d.Content.Find.Execute("44×54=2376", $true, $false, $false, $false, $false, $true, 1, $false, "32×62=1984", 2) | Out-Null
$d.Content.Find.Execute("68×44=2992", $true, $false, $false, $false, $false, $true, 1, $false, "15×11=165", 2) | Out-Null
